$wb = $excel.ActiveWorkbook

# --- Sheet "0D" (sheet1.xml) ---
$ws0D = $wb.Worksheets.Item("0D")
# B2: replace static value with formula "=0.270880692555781*1000"
$ws0D.Range("B2").Formula = "=0.270880692555781*1000"
# Selection moves to B3, tab no longer selected (handled by activating year_Vecteurs at the end)
$ws0D.Range("B3").Select()

# --- Sheet "year_Vecteurs" (sheet4.xml) ---
$wsYV = $wb.Worksheets.Item("year_Vecteurs")
$wsYV.Range("D2").Value = 0.079
$wsYV.Range("C3").Value = 0.187
# D3 formula removed, becomes a plain value
$wsYV.Range("D3").Value = 0.04
$wsYV.Range("C4").Value = 0.272
$wsYV.Range("D4").Value = 0.057
$wsYV.Range("C5").Value = 0.027
$wsYV.Range("C6").Value = 0.3465
$wsYV.Range("D6").Value = 0.0285

# Select G11 on year_Vecteurs and make it the active/selected tab
$wsYV.Activate()
$wsYV.Range("G11").Select()

# --- Window view position ---
$excel.ActiveWindow.WindowState = -4143
$excel.ActiveWindow.Left = 3040
$excel.ActiveWindow.Top = 1000
